$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1209.723
$ws.Range("I15").Value = 1209.723
$ws.Range("K15").Value = 3629.169
$ws.Range("M15").Value = -3460.169
$ws.Range("H17").Value = 3375.0476
$ws.Range("J17").Value = 2828.9473
$ws.Range("L17").Value = 8486.8419
$ws.Range("N17").Value = -8822.8419
$ws.Range("H64").Value = 3099.75
$ws.Range("I64").Value = 2999.5
$ws.Range("K64").Value = 2999.5
$ws.Range("M64").Value = -2751.5
$ws.Range("H67").Value = 3099.75
$ws.Range("I67").Value = 2999.5
$ws.Range("K67").Value = 2999.5
$ws.Range("M67").Value = -2141.5
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H100").Value = 1999.5
$ws.Range("I100").Value = 999.5
$ws.Range("K100").Value = 999.5
$ws.Range("M100").Value = -458.5
$ws.Range("H137").Value = 2018
$ws.Range("I137").Value = 1301
$ws.Range("K137").Value = 3903
$ws.Range("M137").Value = -1353
$ws.Range("H138").Value = 2773.608
$ws.Range("J138").Value = 2696.0908
$ws.Range("L138").Value = 8088.2724
$ws.Range("N138").Value = -18368.2724

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6251.1304
$ws.Range("I61").Value = 7165.2666
$ws.Range("K61").Value = 7165.2666
$ws.Range("M61").Value = -6953.2666
$ws.Range("H63").Value = 1398.8
$ws.Range("I63").Value = 999
$ws.Range("J63").Value = 1998.5
$ws.Range("K63").Value = 999
$ws.Range("L63").Value = 1998.5
$ws.Range("M63").Value = -313
$ws.Range("N63").Value = -3370.5
$ws.Range("H66").Value = 1398.8
$ws.Range("I66").Value = 999
$ws.Range("J66").Value = 1998.5
$ws.Range("K66").Value = 4995
$ws.Range("L66").Value = 9992.5
$ws.Range("M66").Value = -1563
$ws.Range("N66").Value = -16856.5
$ws.Range("H132").Value = 1459.1765
$ws.Range("I132").Value = 1097.28
$ws.Range("J132").Value = 2464.4443
$ws.Range("K132").Value = 3291.84
$ws.Range("L132").Value = 7393.3329
$ws.Range("M132").Value = -761.8400000000001
$ws.Range("N132").Value = -12453.3329
$ws.Range("H136").Value = 6251.1304
$ws.Range("I136").Value = 7165.2666
$ws.Range("K136").Value = 21495.7998
$ws.Range("M136").Value = -18945.7998

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 30000
$ws.Range("J31").Value = 30000
$ws.Range("L31").Value = 30000
$ws.Range("N31").Value = -30504
$ws.Range("H99").Value = 1608
$ws.Range("J99").Value = 1996.5
$ws.Range("L99").Value = 1996.5
$ws.Range("N99").Value = -4992.5
$ws.Range("H107").Value = 1211.9166
$ws.Range("I107").Value = 895.25
$ws.Range("J107").Value = 1370.25
$ws.Range("K107").Value = 895.25
$ws.Range("L107").Value = 1370.25
$ws.Range("M107").Value = 1024.75
$ws.Range("N107").Value = -5210.25
$ws.Range("H134").Value = 6864.773
$ws.Range("I134").Value = 7369.9
$ws.Range("J134").Value = 1813.5
$ws.Range("K134").Value = 22109.7
$ws.Range("L134").Value = 5440.5
$ws.Range("M134").Value = -19574.7
$ws.Range("N134").Value = -10510.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 19000
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 23500
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 23500
$ws.Range("M13").Value = -861
$ws.Range("N13").Value = -23778
$ws.Range("H103").Value = 17499
$ws.Range("I103").Value = 17499
$ws.Range("K103").Value = 17499
$ws.Range("M103").Value = -16327
$ws.Range("H107").Value = 673.3333
$ws.Range("I107").Value = 431.0625
$ws.Range("J107").Value = 1448.6
$ws.Range("K107").Value = 431.0625
$ws.Range("L107").Value = 1448.6
$ws.Range("M107").Value = 1488.9375
$ws.Range("N107").Value = -5288.6
$ws.Range("H132").Value = 1733.1428
$ws.Range("I132").Value = 1029.7826
$ws.Range("K132").Value = 3089.3478
$ws.Range("M132").Value = -559.3478

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 192.25
$ws.Range("I10").Value = 192.25
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 576.75
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -437.75
$ws.Range("N10").ClearContents()
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H113").Value = 5007.5415
$ws.Range("J113").Value = 841.9048
$ws.Range("L113").Value = 2525.7144
$ws.Range("N113").Value = -6865.7144
$ws.Range("H122").Value = 888.625
$ws.Range("I122").Value = 682
$ws.Range("K122").Value = 6138
$ws.Range("M122").Value = -3688
$ws.Range("H132").Value = 2850
$ws.Range("I132").Value = 1200
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 10800
$ws.Range("L132").Value = 40500
$ws.Range("M132").Value = -8270
$ws.Range("N132").Value = -45560
$ws.Range("H133").Value = 4020.9092
$ws.Range("J133").Value = 4775
$ws.Range("L133").Value = 14325
$ws.Range("N133").Value = -24445

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 792504.25
$ws.Range("I22").Value = 1504999.5
$ws.Range("K22").Value = 1504999.5
$ws.Range("M22").Value = -1504470.5
$ws.Range("H80").Value = 2767.0833
$ws.Range("I80").Value = 2519.4
$ws.Range("J80").Value = 4005.5
$ws.Range("K80").Value = 2519.4
$ws.Range("L80").Value = 4005.5
$ws.Range("M80").Value = -1521.4
$ws.Range("N80").Value = -6001.5
$ws.Range("H83").Value = 2767.0833
$ws.Range("I83").Value = 2519.4
$ws.Range("J83").Value = 4005.5
$ws.Range("K83").Value = 12597
$ws.Range("L83").Value = 20027.5
$ws.Range("M83").Value = -7605
$ws.Range("N83").Value = -30011.5
$ws.Range("H126").Value = 3773597.8
$ws.Range("I126").Value = 4633080.5
$ws.Range("K126").Value = 13899241.5
$ws.Range("M126").Value = -13896771.5
$ws.Range("H132").Value = 1375771.6
$ws.Range("I132").Value = 1749718.5
$ws.Range("K132").Value = 5249155.5
$ws.Range("M132").Value = -5246625.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 14000
$ws.Range("J5").Value = 14000
$ws.Range("L5").Value = 14000
$ws.Range("N5").Value = -14226
$ws.Range("H46").Value = 1504.6
$ws.Range("I46").Value = 1064.6923
$ws.Range("K46").Value = 1064.6923
$ws.Range("M46").Value = -876.6922999999999
$ws.Range("H82").Value = 2773.8
$ws.Range("I82").Value = 1997.5
$ws.Range("J82").Value = 3291.3333
$ws.Range("K82").Value = 1997.5
$ws.Range("L82").Value = 3291.3333
$ws.Range("M82").Value = -1636.5
$ws.Range("N82").Value = -4013.3333
$ws.Range("H85").Value = 2773.8
$ws.Range("I85").Value = 1997.5
$ws.Range("J85").Value = 3291.3333
$ws.Range("K85").Value = 1997.5
$ws.Range("L85").Value = 3291.3333
$ws.Range("M85").Value = -749.5
$ws.Range("N85").Value = -5787.3333
$ws.Range("H132").Value = 2813.7273
$ws.Range("I132").Value = 1628.6364
$ws.Range("K132").Value = 4885.9092
$ws.Range("M132").Value = -2355.9092
$ws.Range("H136").Value = 1543.5135
$ws.Range("I136").Value = 1207.3334
$ws.Range("K136").Value = 3622.0002
$ws.Range("M136").Value = -1072.0002

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 699.4583
$ws.Range("I107").Value = 499.11765
$ws.Range("J107").Value = 1186
$ws.Range("K107").Value = 1497.35295
$ws.Range("L107").Value = 3558
$ws.Range("M107").Value = 422.64705
$ws.Range("N107").Value = -7398
$ws.Range("H132").Value = 1705.8846
$ws.Range("I132").Value = 1109.8636
$ws.Range("K132").Value = 3329.5908
$ws.Range("M132").Value = -799.5908
$ws.Range("H136").Value = 30865984
$ws.Range("I136").Value = 55556916
$ws.Range("K136").Value = 166670748
$ws.Range("M136").Value = -166668198
